# Setting up Ubuntu.docx -- apply commit "Add Roku position display and activity check"
# (per the supplied OOXML diff: tidy several bullet runs, rework the JRMC
#  service bullet, turn two URLs into real hyperlinks, add a few new bullets,
#  and split two bullets' text.)
#
# Strategy: use Range.InsertXML() with fully specified <w:p> fragments so the
# resulting markup matches the target exactly (single merged run, no stray
# w:proofErr, explicit rPr where needed). Edits are applied from the bottom
# of the document upward so paragraph indices of not-yet-processed bullets
# stay stable. Hyperlinks are added afterwards with Hyperlinks.Add so Word
# wires up the relationship + Hyperlink character style itself.

$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Paragraph 26: "Install and configure Nginx" URL
#   -> bullet "Install and configure Nginx" stays (handled below, para 25)
#   -> this bullet (the URL) gets a trailing new bullet "Allow port 84 ..."
#      (hyperlink-ification happens afterwards)
# ---------------------------------------------------------------------
$d.Paragraphs(26).Range.InsertXML(
  "<w:p $ns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t>https://learn.microsoft.com/en-us/aspnet/core/host-and-deploy/linux-nginx?view=aspnetcore-6.0&amp;tabs=linux-ubuntu</w:t></w:r></w:p>" +
  "<w:p $ns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t xml:space='preserve'>Allow port 84 through the firewall </w:t></w:r></w:p>"
)

# ---------------------------------------------------------------------
# Paragraph 25: "Install and configure Nginx" -- drop yellow highlight,
# merge into one run, remove proofErr
# ---------------------------------------------------------------------
$d.Paragraphs(25).Range.InsertXML(
  "<w:p $ns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t>Install and configure Nginx</w:t></w:r></w:p>"
)

# ---------------------------------------------------------------------
# Paragraph 24: "Set to start automatically (Can Nginx do this?)"
#   -> "Set to start automatically " + "(in a bash loop)" (two runs)
# ---------------------------------------------------------------------
$d.Paragraphs(24).Range.InsertXML(
  "<w:p $ns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t xml:space='preserve'>Set to start automatically </w:t></w:r><w:r><w:t>(in a bash loop)</w:t></w:r></w:p>"
)

# ---------------------------------------------------------------------
# Paragraph 21: "Copy over Avid5.Net" -- merge into one run
# ---------------------------------------------------------------------
$d.Paragraphs(21).Range.InsertXML(
  "<w:p $ns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t>Copy over Avid5.Net</w:t></w:r></w:p>"
)

# ---------------------------------------------------------------------
# Paragraph 20: "sudo apt-get install -y aspnetcore-runtime-6.0 "
#   -> merge "sudo" + " apt-get install -y aspnetcore-runtime-" into one run,
#      keep "6" / ".0" / " " separate, drop proofErr
# ---------------------------------------------------------------------
$d.Paragraphs(20).Range.InsertXML(
  "<w:p $ns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t>sudo apt-get install -y aspnetcore-runtime-</w:t></w:r><w:r><w:t>6</w:t></w:r><w:r><w:t>.0</w:t></w:r><w:r><w:t xml:space='preserve'> </w:t></w:r></w:p>"
)

# ---------------------------------------------------------------------
# Paragraph 19: "Install .Net 6 runtime" -- merge into one run
# ---------------------------------------------------------------------
$d.Paragraphs(19).Range.InsertXML(
  "<w:p $ns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t>Install .Net 6 runtime</w:t></w:r></w:p>"
)

# ---------------------------------------------------------------------
# Paragraph 18: Samba URL
#   -> stays as its own bullet (hyperlink-ified afterwards), plus two new
#      bullets: "At least media and Avid5 folders" and the firewall note
# ---------------------------------------------------------------------
$d.Paragraphs(18).Range.InsertXML(
  "<w:p $ns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t>https://ubuntu.com/server/docs/samba-file-server</w:t></w:r></w:p>" +
  "<w:p $ns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t>At least media and Avid5 folders</w:t></w:r></w:p>" +
  "<w:p $ns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t xml:space='preserve'>Don" + [char]0x2019 + "t forget the firewall : </w:t></w:r><w:r><w:t>sudo ufw allow samba</w:t></w:r></w:p>"
)

# ---------------------------------------------------------------------
# Paragraph 17: "Install and configure Samba" -- merge into one run
# ---------------------------------------------------------------------
$d.Paragraphs(17).Range.InsertXML(
  "<w:p $ns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t>Install and configure Samba</w:t></w:r></w:p>"
)

# ---------------------------------------------------------------------
# Paragraph 16: "Configure JRMC library, Network, TV etc"
#   -> "Network, " run stays separate; "TV " + "etc" merge into "TV etc"
# ---------------------------------------------------------------------
$d.Paragraphs(16).Range.InsertXML(
  "<w:p $ns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t xml:space='preserve'>Configure JRMC library, </w:t></w:r><w:r><w:t xml:space='preserve'>Network, </w:t></w:r><w:r><w:t>TV etc</w:t></w:r></w:p>"
)

# ---------------------------------------------------------------------
# Paragraph 14: "??? installJRMC for service? " (yellow highlight)
#   -> "installJRMC for service" + " " + "(--service jriver-mediaserver)"
#      (last run gets Verdana/black/sz20/shading, no more highlight)
# ---------------------------------------------------------------------
$d.Paragraphs(14).Range.InsertXML(
  "<w:p $ns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t>installJRMC for service</w:t></w:r><w:r><w:t xml:space='preserve'> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii='Verdana' w:hAnsi='Verdana'/><w:color w:val='000000'/><w:sz w:val='20'/><w:szCs w:val='20'/><w:shd w:val='clear' w:color='auto' w:fill='ECEDF3'/></w:rPr><w:t>(--service jriver-mediaserver)</w:t></w:r></w:p>"
)

# ---------------------------------------------------------------------
# Paragraph 13: "Install and license JRMC" -- merge into one run
# ---------------------------------------------------------------------
$d.Paragraphs(13).Range.InsertXML(
  "<w:p $ns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t>Install and license JRMC</w:t></w:r></w:p>"
)

# ---------------------------------------------------------------------
# Paragraph 12: "Install and login Chrome" -- merge into one run
# ---------------------------------------------------------------------
$d.Paragraphs(12).Range.InsertXML(
  "<w:p $ns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t>Install and login Chrome</w:t></w:r></w:p>"
)

# ---------------------------------------------------------------------
# Paragraph 10: "Mount media disk in /media (en dash) edit /etc/fstab"
#   -- merge into one run
# ---------------------------------------------------------------------
$d.Paragraphs(10).Range.InsertXML(
  "<w:p $ns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t>Mount media disk in /media " + [char]0x2013 + " edit /etc/fstab</w:t></w:r></w:p>"
)

# ---------------------------------------------------------------------
# Paragraph 6: "Install NoMachine for remote access" -- merge into one run
# ---------------------------------------------------------------------
$d.Paragraphs(6).Range.InsertXML(
  "<w:p $ns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t>Install NoMachine for remote access</w:t></w:r></w:p>"
)

# ---------------------------------------------------------------------
# Paragraph 4: "Configure User in settings to automatically login"
#   -- merge into one run
# ---------------------------------------------------------------------
$d.Paragraphs(4).Range.InsertXML(
  "<w:p $ns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t>Configure User in settings to automatically login</w:t></w:r></w:p>"
)

# ---------------------------------------------------------------------
# Paragraph 3: "...image backups and faster root backup" -- merge last
#   two runs ("and faster root" + "backup"), drop proofErr
# ---------------------------------------------------------------------
$d.Paragraphs(3).Range.InsertXML(
  "<w:p $ns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t xml:space='preserve'>Make a small </w:t></w:r><w:r><w:t xml:space='preserve'>(64GB max) </w:t></w:r><w:r><w:t>partition for root, to allow space for image backups</w:t></w:r><w:r><w:t xml:space='preserve'> and faster root backup</w:t></w:r></w:p>"
)

# ---------------------------------------------------------------------
# Turn the two bare URLs into real hyperlinks (adds the relationship +
# applies the built-in Hyperlink character style), matching the diff's
# <w:hyperlink> wrapped runs.
# ---------------------------------------------------------------------
$sambaUrl = "https://ubuntu.com/server/docs/samba-file-server"
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq $sambaUrl) {
        $r = $p.Range
        $r.MoveEnd(1, -1)
        $r.Delete()
        $d.Hyperlinks.Add($r, $sambaUrl, "", "", $sambaUrl) | Out-Null
        break
    }
}

$nginxUrl = "https://learn.microsoft.com/en-us/aspnet/core/host-and-deploy/linux-nginx?view=aspnetcore-6.0&tabs=linux-ubuntu"
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq $nginxUrl) {
        $r = $p.Range
        $r.MoveEnd(1, -1)
        $r.Delete()
        $d.Hyperlinks.Add($r, $nginxUrl, "", "", $nginxUrl) | Out-Null
        break
    }
}

Write-Output "done"
